$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 111
$ws.Range("H111").Value = 3932.353
$ws.Range("I111").Value = 3487.5
$ws.Range("J111").Value = 5000
$ws.Range("K111").Value = 10462.5
$ws.Range("L111").Value = 15000
$ws.Range("M111").Value = -7395.5
$ws.Range("N111").Value = -21134
# Row 112
$ws.Range("H112").Value = 1494.1666
$ws.Range("J112").Value = 1494.1666
$ws.Range("L112").Value = 4482.4998
$ws.Range("N112").Value = -6698.4998
# Row 138
$ws.Range("H138").Value = 2224668.2
$ws.Range("I138").Value = 1024.1613
$ws.Range("J138").Value = 3393023.5
$ws.Range("K138").Value = 3072.4839
$ws.Range("L138").Value = 10179070.5
$ws.Range("M138").Value = 2067.5161
$ws.Range("N138").Value = -10189350.5

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 18531.518
$ws.Range("I74").Value = 24234.744
$ws.Range("J74").Value = 2182.2666
$ws.Range("K74").Value = 24234.744
$ws.Range("L74").Value = 2182.2666
$ws.Range("M74").Value = -23360.744
$ws.Range("N74").Value = -3930.2666
# Row 77
$ws.Range("H77").Value = 18531.518
$ws.Range("I77").Value = 24234.744
$ws.Range("J77").Value = 2182.2666
$ws.Range("K77").Value = 121173.72
$ws.Range("L77").Value = 10911.333
$ws.Range("M77").Value = -116805.72
$ws.Range("N77").Value = -19647.333
# Row 118
$ws.Range("H118").Value = 28641.223
$ws.Range("J118").Value = 30050
$ws.Range("L118").Value = 30050
$ws.Range("N118").Value = -33364
# Row 132
$ws.Range("H132").Value = 1858.1154
$ws.Range("I132").Value = 1914.591
$ws.Range("J132").Value = 1547.5
$ws.Range("K132").Value = 5743.772999999999
$ws.Range("L132").Value = 4642.5
$ws.Range("M132").Value = -3213.772999999999
$ws.Range("N132").Value = -9702.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2207.818
$ws.Range("I105").Value = 2098
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 2098
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -351
$ws.Range("N105").Value = -5894

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 1003010.5
$ws.Range("I132").Value = 2011.0312
$ws.Range("J132").Value = 7409407
$ws.Range("K132").Value = 6033.0936
$ws.Range("L132").Value = 22228221
$ws.Range("M132").Value = -3503.0936
$ws.Range("N132").Value = -22233281
# Row 141
$ws.Range("H141").Value = 69936.664
$ws.Range("J141").Value = 77864.8
$ws.Range("L141").Value = 77864.8
$ws.Range("N141").Value = -88224.8

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 568
$ws.Range("I16").Value = 490
$ws.Range("J16").Value = 802
$ws.Range("K16").Value = 1470
$ws.Range("L16").Value = 2406
$ws.Range("M16").Value = -1297
$ws.Range("N16").Value = -2752
# Row 98
$ws.Range("H98").Value = 531.2143
$ws.Range("I98").Value = 443.125
$ws.Range("J98").Value = 648.6667
$ws.Range("K98").Value = 1329.375
$ws.Range("L98").Value = 1946.0001
$ws.Range("M98").Value = 168.625
$ws.Range("N98").Value = -4942.0001
# Row 107
$ws.Range("H107").Value = 333743.1
$ws.Range("I107").Value = 334.33334
$ws.Range("J107").Value = 370788.53
$ws.Range("K107").Value = 1003.00002
$ws.Range("L107").Value = 1112365.59
$ws.Range("M107").Value = 916.9999799999999
$ws.Range("N107").Value = -1116205.59
# Row 110
$ws.Range("H110").Value = 1740.875
$ws.Range("I110").Value = 1740.875
$ws.Range("K110").Value = 5222.625
$ws.Range("M110").Value = -1132.625
# Row 120
$ws.Range("H120").Value = 8447.777
$ws.Range("I120").Value = 7432.857
$ws.Range("K120").Value = 22298.571
$ws.Range("M120").Value = -17460.571
# Row 131
$ws.Range("H131").Value = 919.72
$ws.Range("J131").Value = 920.7835
$ws.Range("L131").Value = 2762.3505
$ws.Range("N131").Value = -12842.3505

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Range("H47").Value = 11000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 11000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 11000
$ws.Range("M47").Value = ""
$ws.Range("N47").Value = -12136
# Row 102
$ws.Range("H102").Value = 8334777
$ws.Range("I102").Value = 10417844
$ws.Range("J102").Value = 2512
$ws.Range("K102").Value = 10417844
$ws.Range("L102").Value = 2512
$ws.Range("M102").Value = -10416222
$ws.Range("N102").Value = -5756
# Row 122
$ws.Range("H122").Value = 148222.14
$ws.Range("I122").Value = 170333.33
$ws.Range("J122").Value = 15555
$ws.Range("K122").Value = 510999.99
$ws.Range("L122").Value = 46665
$ws.Range("M122").Value = -508549.99
$ws.Range("N122").Value = -51565
# Row 124
$ws.Range("H124").Value = 25250
$ws.Range("J124").Value = 25250
$ws.Range("L124").Value = 25250
$ws.Range("N124").Value = -35070
# Row 132
$ws.Range("H132").Value = 2634412.8
$ws.Range("I132").Value = 2538.394
$ws.Range("J132").Value = 20004782
$ws.Range("K132").Value = 7615.181999999999
$ws.Range("L132").Value = 60014346
$ws.Range("M132").Value = -5085.181999999999
$ws.Range("N132").Value = -60019406

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2900
$ws.Range("I7").Value = 2875
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2875
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2763
$ws.Range("N7").Value = -3224
# Row 61
$ws.Range("H61").Value = 2512.5
$ws.Range("I61").Value = 2050
$ws.Range("J61").Value = 2975
$ws.Range("K61").Value = 2050
$ws.Range("L61").Value = 2975
$ws.Range("M61").Value = -1848
$ws.Range("N61").Value = -3379
# Row 68
$ws.Range("H68").Value = 9971.357
$ws.Range("I68").Value = 16357.714
$ws.Range("J68").Value = 3585
$ws.Range("K68").Value = 16357.714
$ws.Range("L68").Value = 3585
$ws.Range("M68").Value = -15608.714
$ws.Range("N68").Value = -5083
# Row 71
$ws.Range("H71").Value = 9971.357
$ws.Range("I71").Value = 16357.714
$ws.Range("J71").Value = 3585
$ws.Range("K71").Value = 81788.57
$ws.Range("L71").Value = 17925
$ws.Range("M71").Value = -78044.57
$ws.Range("N71").Value = -25413
# Row 113
$ws.Range("H113").Value = 2512.5
$ws.Range("I113").Value = 2050
$ws.Range("J113").Value = 2975
$ws.Range("K113").Value = 2050
$ws.Range("L113").Value = 2975
$ws.Range("M113").Value = 120
$ws.Range("N113").Value = -7315
# Row 126
$ws.Range("H126").Value = 2900
$ws.Range("I126").Value = 2875
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8625
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6155
$ws.Range("N126").Value = -13940
# Row 132
$ws.Range("H132").Value = 2916.4443
$ws.Range("I132").Value = 2916.4443
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8749.332900000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6219.332900000001
$ws.Range("N132").Value = ""

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 15000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = ""
# Row 109
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774
# Row 132
$ws.Range("H132").Value = 2517.6123
$ws.Range("I132").Value = 3038.4333
$ws.Range("J132").Value = 1695.2632
$ws.Range("K132").Value = 9115.2999
$ws.Range("L132").Value = 5085.7896
$ws.Range("M132").Value = -6585.2999
$ws.Range("N132").Value = -10145.7896
# Row 133
$ws.Range("H133").Value = 30943
$ws.Range("J133").Value = 30943
$ws.Range("L133").Value = 30943
$ws.Range("N133").Value = -41063
